$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tournament")
$ws.Activate()

# Insert two new rows above the old row 2 ("name" / "location" / ...),
# shifting the existing tournament-level rows down by two. These new
# rows will carry the competition-key / host-key lookup rows.
$ws.Range("A2:A3").EntireRow.Insert()

$ws.Range("A2").Value = "competition-key"
$ws.Range("A3").Value = "host-key"
$ws.Range("B3").Value = "usa"
$ws.Range("B2").Value = "mens-club-world-cup"

# Append the new venue-key rows after the previously-last row (old row 15,
# now row 17).
$venueKeys = @(
    @("venue-key.1", "us-atlanta-ga"),
    @("venue-key.2", "us-charlotte-nc"),
    @("venue-key.3", "us-cincinnati-oh"),
    @("venue-key.4", "us-pasadena-ca"),
    @("venue-key.5", "us-miami-fl"),
    @("venue-key.6", "us-nashville-tn"),
    @("venue-key.7", "us-east-rutherford-nj"),
    @("venue-key.8", "us-orlando-fl"),
    @("venue-key.9", "us-philadelphia-pa"),
    @("venue-key.10", "us-seattle-wa"),
    @("venue-key.11", "us-washington-dc")
)

$startRow = 18
for ($i = 0; $i -lt $venueKeys.Count; $i++) {
    $r = $startRow + $i
    $pair = $venueKeys[$i]
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
}

# Resize the "tournament" table (and its autofilter) to cover the new extent.
$lo = $ws.ListObjects.Item("tournament")
$lo.Resize($ws.Range("A1:I28"))

# Restore the active selection to the newly-inserted rows.
$ws.Range("A2:XFD3").Select()
